# Add our VIVA problems here
# - Turn the old "5. Are we ready?" / "????" rows into new problem text,
#   and move the "Are we ready?" / "????" items further down the list
#   as new rows (10 and 11), leaving rows 8-9 blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new rows first so the shared-string table order matches
# the target (the moved-down items come right after the existing ones).
$ws.Range("A10").Value = "Are we ready?"

# Replace rows 6 and 7 with the new VIVA problem text.
$ws.Range("A6").Value = "5. What about the fonts for slides? Is it TimesNewRoman?"
$ws.Range("A7").Value = "6. Reviews for 1st slide"

$ws.Range("A11").Value = "????"

# Match the author's recorded selection after the edit.
$ws.Range("A15").Select() | Out-Null
